$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.142.87"
$ws.Range("E2").Value = "  +3.46%  "
$ws.Range("D3").Value = "1.604.23"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Formula = "'212.87"
$ws.Range("E5").Value = "  +3.13%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Formula = "'0.487"
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("D12").Value = "1.827.32"
$ws.Range("E12").Value = "  +3.48%  "
$ws.Range("D13").Value = "1.611.22"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "26.142.79"
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  +10.29%  "
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("E24").Value = "  +15.24%  "
$ws.Range("D25").Formula = "'141.61"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").Formula = "'0.0471"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").Formula = "'0.0165"
$ws.Range("E36").Value = "  +11.11%  "
$ws.Range("D37").Value = "1.123.83"
$ws.Range("E37").Value = "  +4.03%  "
$ws.Range("E39").Value = "  +3.10%  "
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Formula = "'0.493"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").Formula = "'5.16"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("D44").Value = "1.739.62"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("D45").Formula = "'92.78"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("D47").Formula = "'53.54"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -16.61%  "
